$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Insert a new blank column before column N (14th column), shifting the
# existing "Late" / "Outstanding" (and the blank spacer) columns one to
# the right - this is the "Variable Instalments" column added for the
# Loan RBI change.
$ws.Columns("N").Insert()

# The newly inserted column picks up the width of the column to its
# left (M - "In Advance"), matching Excel's default "insert column"
# formatting behaviour.
$ws.Columns("N").ColumnWidth = $ws.Columns("M").ColumnWidth

# Restore the active selection on the sheet (it shifted from C8 to P8
# because of the inserted column).
$ws.Range("P8").Select()
